# Update Name of Algo
# Apply updated imputed values to the RandomForest result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.6515
$ws.Range("B3").Value = 5.979700000000005
$ws.Range("B14").Value = 5.346200000000002
$ws.Range("B16").Value = 6.650499999999997
$ws.Range("D18").Value = -8.855299999999998
$ws.Range("B21").Value = 9.265200000000004
$ws.Range("B23").Value = 9.016500000000006
$ws.Range("D24").Value = -7.269999999999999
$ws.Range("B25").Value = 5.5582
$ws.Range("D25").Value = -8.063099999999999
$ws.Range("B26").Value = 6.7861
$ws.Range("D27").Value = -8.891800000000002
$ws.Range("B29").Value = 5.086800000000002
$ws.Range("D30").Value = -7.4809
$ws.Range("D31").Value = -8.272600000000001
$ws.Range("D39").Value = -8.163299999999998
$ws.Range("B40").Value = 8.959699999999996
$ws.Range("D42").Value = -8.257100000000001
$ws.Range("D48").Value = -7.422299999999999
$ws.Range("D51").Value = -7.794899999999996
$ws.Range("D52").Value = -7.998200000000002
$ws.Range("B53").Value = 5.2001
$ws.Range("D55").Value = -8.8317
$ws.Range("D56").Value = -7.884099999999997
$ws.Range("B57").Value = 4.939799999999998
$ws.Range("D57").Value = -8.031899999999997
$ws.Range("B59").Value = 6.518900000000002
$ws.Range("D60").Value = -7.901199999999999
$ws.Range("B65").Value = 5.707300000000001
$ws.Range("B69").Value = 5.377699999999995
$ws.Range("D73").Value = -8.075999999999995
$ws.Range("D74").Value = -8.138100000000007
$ws.Range("B79").Value = 9.431700000000005
$ws.Range("B83").Value = 4.938099999999996
$ws.Range("D89").Value = -5.879400000000001
$ws.Range("D90").Value = -8.195200000000005
$ws.Range("B91").Value = 4.973600000000001
$ws.Range("D92").Value = -5.782800000000002
$ws.Range("B93").Value = 5.805400000000001
$ws.Range("B100").Value = 5.050999999999997

$wb.Save()
